$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.5123330574880072
$ws.Cells.Item(2, 3).Value = 0.1191298221600121
$ws.Cells.Item(2, 4).Value = 0.04629170208831823
$ws.Cells.Item(2, 5).Value = 0.1009989519192374
$ws.Cells.Item(2, 6).Value = 0.9232754284974618
$ws.Cells.Item(2, 7).Value = 0
$ws.Cells.Item(2, 8).Value = 0.07973214163530429
$ws.Cells.Item(2, 9).Value = 0.8623522309721103
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 11).Value = 0.3285408004581996
$ws.Cells.Item(2, 12).Value = 0.2043011126019394
$ws.Cells.Item(2, 13).Value = 0
$ws.Cells.Item(2, 14).Value = 0
$ws.Cells.Item(2, 15).Value = 3.289596083056225

$ws.Cells.Item(3, 2).Value = 0.4690635668912932
$ws.Cells.Item(3, 3).Value = 0.1177010559582499
$ws.Cells.Item(3, 4).Value = 0.04408252612715557
$ws.Cells.Item(3, 5).Value = 0.1005822379081991
$ws.Cells.Item(3, 6).Value = 0.9252923456370752
$ws.Cells.Item(3, 7).Value = 0
$ws.Cells.Item(3, 8).Value = 0.07973214163530429
$ws.Cells.Item(3, 9).Value = 0.869696731390853
$ws.Cells.Item(3, 10).Value = 0
$ws.Cells.Item(3, 11).Value = 0.2894826279240021
$ws.Cells.Item(3, 12).Value = 0.1968817526098974
$ws.Cells.Item(3, 13).Value = 0
$ws.Cells.Item(3, 14).Value = 0
$ws.Cells.Item(3, 15).Value = 3.310028172141486

$ws.Cells.Item(4, 2).Value = 0.4425834031724207
$ws.Cells.Item(4, 3).Value = 0.1168192702539415
$ws.Cells.Item(4, 4).Value = 0.04271260762028817
$ws.Cells.Item(4, 5).Value = 0.1003781762352318
$ws.Cells.Item(4, 6).Value = 0.9270721256038712
$ws.Cells.Item(4, 7).Value = 0
$ws.Cells.Item(4, 8).Value = 0.07973214163530429
$ws.Cells.Item(4, 9).Value = 0.8746426740282693
$ws.Cells.Item(4, 10).Value = 0
$ws.Cells.Item(4, 11).Value = 0.2654701821502954
$ws.Cells.Item(4, 12).Value = 0.1924331771722478
$ws.Cells.Item(4, 13).Value = 0
$ws.Cells.Item(4, 14).Value = 0
$ws.Cells.Item(4, 15).Value = 3.324370843595617

$ws.Cells.Item(5, 2).Value = 0.4318152270989515
$ws.Cells.Item(5, 3).Value = 0.1164588272108276
$ws.Cells.Item(5, 4).Value = 0.04215099352329332
$ws.Cells.Item(5, 5).Value = 0.1003080642353389
$ws.Cells.Item(5, 6).Value = 0.9279335626859222
$ws.Cells.Item(5, 7).Value = 0
$ws.Cells.Item(5, 8).Value = 0.07973214163530429
$ws.Cells.Item(5, 9).Value = 0.876767926181973
$ws.Cells.Item(5, 10).Value = 0
$ws.Cells.Item(5, 11).Value = 0.2556778000440119
$ws.Cells.Item(5, 12).Value = 0.1906473205182522
$ws.Cells.Item(5, 13).Value = 0
$ws.Cells.Item(5, 14).Value = 0
$ws.Cells.Item(5, 15).Value = 3.330667545683099

$ws.Cells.Item(6, 2).Value = 0.4300285719298813
$ws.Cells.Item(6, 3).Value = 0.116398909806172
$ws.Cells.Item(6, 4).Value = 0.04205753572987447
$ws.Cells.Item(6, 5).Value = 0.1002972106127373
$ws.Cells.Item(6, 6).Value = 0.9280848283095935
$ws.Cells.Item(6, 7).Value = 0
$ws.Cells.Item(6, 8).Value = 0.07973214163530429
$ws.Cells.Item(6, 9).Value = 0.8771274509294074
$ws.Cells.Item(6, 10).Value = 0
$ws.Cells.Item(6, 11).Value = 0.2540513695667386
$ws.Cells.Item(6, 12).Value = 0.1903524119767468
$ws.Cells.Item(6, 13).Value = 0
$ws.Cells.Item(6, 14).Value = 0
$ws.Cells.Item(6, 15).Value = 3.331740404910875

$ws.Cells.Item(7, 2).Value = 0.4424380869089646
$ws.Cells.Item(7, 3).Value = 0.1168144136428921
$ws.Cells.Item(7, 4).Value = 0.0427050470686936
$ws.Cells.Item(7, 5).Value = 0.1003771778376183
$ws.Cells.Item(7, 6).Value = 0.9270831918951643
$ws.Cells.Item(7, 7).Value = 0
$ws.Cells.Item(7, 8).Value = 0.07973214163530429
$ws.Cells.Item(7, 9).Value = 0.8746708915793562
$ws.Cells.Item(7, 10).Value = 0
$ws.Cells.Item(7, 11).Value = 0.2653381467112013
$ws.Cells.Item(7, 12).Value = 0.1924089831579181
$ws.Cells.Item(7, 13).Value = 0
$ws.Cells.Item(7, 14).Value = 0
$ws.Cells.Item(7, 15).Value = 3.324453933202548

$ws.Cells.Item(8, 2).Value = 0.4973959915840851
$ws.Cells.Item(8, 3).Value = 0.118638137342316
$ws.Cells.Item(8, 4).Value = 0.04553279084195339
$ws.Cells.Item(8, 5).Value = 0.100844528149473
$ws.Cells.Item(8, 6).Value = 0.9238585112105255
$ws.Cells.Item(8, 7).Value = 0
$ws.Cells.Item(8, 8).Value = 0.07973214163530429
$ws.Cells.Item(8, 9).Value = 0.8647940459803252
$ws.Cells.Item(8, 10).Value = 0
$ws.Cells.Item(8, 11).Value = 0.3150802177987657
$ws.Cells.Item(8, 12).Value = 0.2017207683886681
$ws.Cells.Item(8, 13).Value = 0
$ws.Cells.Item(8, 14).Value = 0
$ws.Cells.Item(8, 15).Value = 3.296268040646083

$ws.Cells.Item(9, 2).Value = 0.6058367801736892
$ws.Cells.Item(9, 3).Value = 0.1221775609745208
$ws.Cells.Item(9, 4).Value = 0.05097010280109515
$ws.Cells.Item(9, 5).Value = 0.102171380338536
$ws.Cells.Item(9, 6).Value = 0.9218304246864619
$ws.Cells.Item(9, 7).Value = 0
$ws.Cells.Item(9, 8).Value = 0.07973214163530429
$ws.Cells.Item(9, 9).Value = 0.8488880017014075
$ws.Cells.Item(9, 10).Value = 0
$ws.Cells.Item(9, 11).Value = 0.4123613429555633
$ws.Cells.Item(9, 12).Value = 0.2208272049869606
$ws.Cells.Item(9, 13).Value = 0
$ws.Cells.Item(9, 14).Value = 0
$ws.Cells.Item(9, 15).Value = 3.255257870112388

$ws.Cells.Item(10, 2).Value = 0.6858895072661539
$ws.Cells.Item(10, 3).Value = 0.1247543254561378
$ws.Cells.Item(10, 4).Value = 0.05489818443147954
$ws.Cells.Item(10, 5).Value = 0.1033957758699131
$ws.Cells.Item(10, 6).Value = 0.9229603830096877
$ws.Cells.Item(10, 7).Value = 0
$ws.Cells.Item(10, 8).Value = 0.07973214163530429
$ws.Cells.Item(10, 9).Value = 0.8393126654137504
$ws.Cells.Item(10, 10).Value = 0
$ws.Cells.Item(10, 11).Value = 0.4836531176303822
$ws.Cells.Item(10, 12).Value = 0.2353788483306971
$ws.Cells.Item(10, 13).Value = 0
$ws.Cells.Item(10, 14).Value = 0
$ws.Cells.Item(10, 15).Value = 3.233827271308002

$ws.Cells.Item(11, 2).Value = 0.7223853077837248
$ws.Cells.Item(11, 3).Value = 0.1259211991945932
$ws.Cells.Item(11, 4).Value = 0.05667050661892858
$ws.Cells.Item(11, 5).Value = 0.1040068576514663
$ws.Cells.Item(11, 6).Value = 0.9240437183541275
$ws.Cells.Item(11, 7).Value = 0
$ws.Cells.Item(11, 8).Value = 0.07973214163530429
$ws.Cells.Item(11, 9).Value = 0.8354149840221368
$ws.Cells.Item(11, 10).Value = 0
$ws.Cells.Item(11, 11).Value = 0.5160424854635153
$ws.Cells.Item(11, 12).Value = 0.2421101631598077
$ws.Cells.Item(11, 13).Value = 0
$ws.Cells.Item(11, 14).Value = 0
$ws.Cells.Item(11, 15).Value = 3.225968220368202

$ws.Cells.Item(12, 2).Value = 0.7362161328603065
$ws.Cells.Item(12, 3).Value = 0.1263622761533014
$ws.Cells.Item(12, 4).Value = 0.05733951892217704
$ws.Cells.Item(12, 5).Value = 0.1042460226719868
$ws.Cells.Item(12, 6).Value = 0.9245358177941725
$ws.Cells.Item(12, 7).Value = 0
$ws.Cells.Item(12, 8).Value = 0.07973214163530429
$ws.Cells.Item(12, 9).Value = 0.8340049254962452
$ws.Cells.Item(12, 10).Value = 0
$ws.Cells.Item(12, 11).Value = 0.5283010296032273
$ws.Cells.Item(12, 12).Value = 0.2446751393163282
$ws.Cells.Item(12, 13).Value = 0
$ws.Cells.Item(12, 14).Value = 0
$ws.Cells.Item(12, 15).Value = 3.223264026484685

$ws.Cells.Item(13, 2).Value = 0.7332369553317335
$ws.Cells.Item(13, 3).Value = 0.126267318062915
$ws.Cells.Item(13, 4).Value = 0.05719553026893465
$ws.Cells.Item(13, 5).Value = 0.1041941692877622
$ws.Cells.Item(13, 6).Value = 0.9244261945103389
$ws.Cells.Item(13, 7).Value = 0
$ws.Cells.Item(13, 8).Value = 0.07973214163530429
$ws.Cells.Item(13, 9).Value = 0.8343056751670019
$ws.Cells.Item(13, 10).Value = 0
$ws.Cells.Item(13, 11).Value = 0.5256612364967737
$ws.Cells.Item(13, 12).Value = 0.2441220168104934
$ws.Cells.Item(13, 13).Value = 0
$ws.Cells.Item(13, 14).Value = 0
$ws.Cells.Item(13, 15).Value = 3.223834330453457

$ws.Cells.Item(14, 2).Value = 0.7235229688529898
$ws.Cells.Item(14, 3).Value = 0.1259575029172169
$ws.Cells.Item(14, 4).Value = 0.05672558948079853
$ws.Cells.Item(14, 5).Value = 0.1040263784486939
$ws.Cells.Item(14, 6).Value = 0.9240825630160998
$ws.Cells.Item(14, 7).Value = 0
$ws.Cells.Item(14, 8).Value = 0.07973214163530429
$ws.Cells.Item(14, 9).Value = 0.8352976566218757
$ws.Cells.Item(14, 10).Value = 0
$ws.Cells.Item(14, 11).Value = 0.5170511392269646
$ws.Cells.Item(14, 12).Value = 0.2423208656733919
$ws.Cells.Item(14, 13).Value = 0
$ws.Cells.Item(14, 14).Value = 0
$ws.Cells.Item(14, 15).Value = 3.225740295690258

$ws.Cells.Item(15, 2).Value = 0.7175742369852003
$ws.Cells.Item(15, 3).Value = 0.1257676280814835
$ws.Cells.Item(15, 4).Value = 0.05643745918341381
$ws.Cells.Item(15, 5).Value = 0.103924611980851
$ws.Cells.Item(15, 6).Value = 0.9238827398001987
$ws.Cells.Item(15, 7).Value = 0
$ws.Cells.Item(15, 8).Value = 0.07973214163530429
$ws.Cells.Item(15, 9).Value = 0.8359138578321321
$ws.Cells.Item(15, 10).Value = 0
$ws.Cells.Item(15, 11).Value = 0.5117763272508569
$ws.Cells.Item(15, 12).Value = 0.2412196861462945
$ws.Cells.Item(15, 13).Value = 0
$ws.Cells.Item(15, 14).Value = 0
$ws.Cells.Item(15, 15).Value = 3.226943161381428

$ws.Cells.Item(16, 2).Value = 0.6835059536811343
$ws.Cells.Item(16, 3).Value = 0.1246779581637725
$ws.Cells.Item(16, 4).Value = 0.05478206295052956
$ws.Cells.Item(16, 5).Value = 0.1033569273259651
$ws.Cells.Item(16, 6).Value = 0.9229010403202125
$ws.Cells.Item(16, 7).Value = 0
$ws.Cells.Item(16, 8).Value = 0.07973214163530429
$ws.Cells.Item(16, 9).Value = 0.8395766098472919
$ws.Cells.Item(16, 10).Value = 0
$ws.Cells.Item(16, 11).Value = 0.4815355056018973
$ws.Cells.Item(16, 12).Value = 0.2349411811617585
$ws.Cells.Item(16, 13).Value = 0
$ws.Cells.Item(16, 14).Value = 0
$ws.Cells.Item(16, 15).Value = 3.23437891586687

$ws.Cells.Item(17, 2).Value = 0.6626259621279473
$ws.Cells.Item(17, 3).Value = 0.1240081001759208
$ws.Cells.Item(17, 4).Value = 0.05376277428631937
$ws.Cells.Item(17, 5).Value = 0.1030225144909025
$ws.Cells.Item(17, 6).Value = 0.9224446124916739
$ws.Cells.Item(17, 7).Value = 0
$ws.Cells.Item(17, 8).Value = 0.07973214163530429
$ws.Cells.Item(17, 9).Value = 0.8419409580604338
$ws.Cells.Item(17, 10).Value = 0
$ws.Cells.Item(17, 11).Value = 0.462972638518778
$ws.Cells.Item(17, 12).Value = 0.2311180674700637
$ws.Cells.Item(17, 13).Value = 0
$ws.Cells.Item(17, 14).Value = 0
$ws.Cells.Item(17, 15).Value = 3.239424615939384

$ws.Cells.Item(18, 2).Value = 0.6506238634293311
$ws.Cells.Item(18, 3).Value = 0.1236223175777624
$ws.Cells.Item(18, 4).Value = 0.05317513572825305
$ws.Cells.Item(18, 5).Value = 0.102835262104719
$ws.Cells.Item(18, 6).Value = 0.9222356738047139
$ws.Cells.Item(18, 7).Value = 0
$ws.Cells.Item(18, 8).Value = 0.07973214163530429
$ws.Cells.Item(18, 9).Value = 0.8433439940478351
$ws.Cells.Item(18, 10).Value = 0
$ws.Cells.Item(18, 11).Value = 0.4522918881025078
$ws.Cells.Item(18, 12).Value = 0.228929632602302
$ws.Cells.Item(18, 13).Value = 0
$ws.Cells.Item(18, 14).Value = 0
$ws.Cells.Item(18, 15).Value = 3.242504643479919

$ws.Cells.Item(19, 2).Value = 0.6465614738891929
$ws.Cells.Item(19, 3).Value = 0.1234916135818622
$ws.Cells.Item(19, 4).Value = 0.05297593693433811
$ws.Cells.Item(19, 5).Value = 0.1027727369326215
$ws.Cells.Item(19, 6).Value = 0.9221741343498522
$ws.Cells.Item(19, 7).Value = 0
$ws.Cells.Item(19, 8).Value = 0.07973214163530429
$ws.Cells.Item(19, 9).Value = 0.8438264440854084
$ws.Cells.Item(19, 10).Value = 0
$ws.Cells.Item(19, 11).Value = 0.4486749220789648
$ws.Cells.Item(19, 12).Value = 0.2281904755097344
$ws.Cells.Item(19, 13).Value = 0
$ws.Cells.Item(19, 14).Value = 0
$ws.Cells.Item(19, 15).Value = 3.243578034833916

$ws.Cells.Item(20, 2).Value = 0.6648479000897112
$ws.Cells.Item(20, 3).Value = 0.1240794594570502
$ws.Cells.Item(20, 4).Value = 0.05387142138765455
$ws.Cells.Item(20, 5).Value = 0.1030575863258782
$ws.Cells.Item(20, 6).Value = 0.9224876540218858
$ws.Cells.Item(20, 7).Value = 0
$ws.Cells.Item(20, 8).Value = 0.07973214163530429
$ws.Cells.Item(20, 9).Value = 0.8416848056750545
$ws.Cells.Item(20, 10).Value = 0
$ws.Cells.Item(20, 11).Value = 0.4649490935107963
$ws.Cells.Item(20, 12).Value = 0.2315239564140512
$ws.Cells.Item(20, 13).Value = 0
$ws.Cells.Item(20, 14).Value = 0
$ws.Cells.Item(20, 15).Value = 3.238869082479852

$ws.Cells.Item(21, 2).Value = 0.7263759189740426
$ws.Cells.Item(21, 3).Value = 0.1260485248847729
$ws.Cells.Item(21, 4).Value = 0.05686368047116019
$ws.Cells.Item(21, 5).Value = 0.1040754521785914
$ws.Cells.Item(21, 6).Value = 0.9241812741798512
$ws.Cells.Item(21, 7).Value = 0
$ws.Cells.Item(21, 8).Value = 0.07973214163530429
$ws.Cells.Item(21, 9).Value = 0.8350044987010961
$ws.Cells.Item(21, 10).Value = 0
$ws.Cells.Item(21, 11).Value = 0.5195803177526557
$ws.Cells.Item(21, 12).Value = 0.2428494746227443
$ws.Cells.Item(21, 13).Value = 0
$ws.Cells.Item(21, 14).Value = 0
$ws.Cells.Item(21, 15).Value = 3.225173088543499

$ws.Cells.Item(22, 2).Value = 0.7666498039228031
$ws.Cells.Item(22, 3).Value = 0.1273307884709709
$ws.Cells.Item(22, 4).Value = 0.05880687897974468
$ws.Cells.Item(22, 5).Value = 0.1047859187120821
$ws.Cells.Item(22, 6).Value = 0.9257653160497625
$ws.Cells.Item(22, 7).Value = 0
$ws.Cells.Item(22, 8).Value = 0.07973214163530429
$ws.Cells.Item(22, 9).Value = 0.8310227090928635
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 11).Value = 0.5552461858651441
$ws.Cells.Item(22, 12).Value = 0.2503444186268808
$ws.Cells.Item(22, 13).Value = 0
$ws.Cells.Item(22, 14).Value = 0
$ws.Cells.Item(22, 15).Value = 3.21780660806175

$ws.Cells.Item(23, 2).Value = 0.7451494730673573
$ws.Cells.Item(23, 3).Value = 0.1266468539510228
$ws.Cells.Item(23, 4).Value = 0.05777090356419734
$ws.Cells.Item(23, 5).Value = 0.1044025964408704
$ws.Cells.Item(23, 6).Value = 0.9248762237614656
$ws.Cells.Item(23, 7).Value = 0
$ws.Cells.Item(23, 8).Value = 0.07973214163530429
$ws.Cells.Item(23, 9).Value = 0.833112705668519
$ws.Cells.Item(23, 10).Value = 0
$ws.Cells.Item(23, 11).Value = 0.5362143955424301
$ws.Cells.Item(23, 12).Value = 0.246335740164227
$ws.Cells.Item(23, 13).Value = 0
$ws.Cells.Item(23, 14).Value = 0
$ws.Cells.Item(23, 15).Value = 3.221593209157021

$ws.Cells.Item(24, 2).Value = 0.6638433553811183
$ws.Cells.Item(24, 3).Value = 0.1240472000192341
$ws.Cells.Item(24, 4).Value = 0.05382230712677938
$ws.Cells.Item(24, 5).Value = 0.103041714755026
$ws.Cells.Item(24, 6).Value = 0.9224680283973186
$ws.Cells.Item(24, 7).Value = 0
$ws.Cells.Item(24, 8).Value = 0.07973214163530429
$ws.Cells.Item(24, 9).Value = 0.8418004759198041
$ws.Cells.Item(24, 10).Value = 0
$ws.Cells.Item(24, 11).Value = 0.4640555652966327
$ws.Cells.Item(24, 12).Value = 0.2313404243500798
$ws.Cells.Item(24, 13).Value = 0
$ws.Cells.Item(24, 14).Value = 0
$ws.Cells.Item(24, 15).Value = 3.239119681004439

$ws.Cells.Item(25, 2).Value = 0.5764318291403754
$ws.Cells.Item(25, 3).Value = 0.1212241234455149
$ws.Cells.Item(25, 4).Value = 0.04951081336815832
$ws.Cells.Item(25, 5).Value = 0.1017685615904185
$ws.Cells.Item(25, 6).Value = 0.9219190703958375
$ws.Cells.Item(25, 7).Value = 0
$ws.Cells.Item(25, 8).Value = 0.07973214163530429
$ws.Cells.Item(25, 9).Value = 0.8528203267250447
$ws.Cells.Item(25, 10).Value = 0
$ws.Cells.Item(25, 11).Value = 0.3860744843158841
$ws.Cells.Item(25, 12).Value = 0.2412196861462945
$ws.Cells.Item(25, 13).Value = 0
$ws.Cells.Item(25, 14).Value = 0
$ws.Cells.Item(25, 15).Value = 3.264824820660806

Write-Output "updated pl_mw data for 380 kV case"